# Retrained the model. Addressing issue #15
# Update yolov3_summary sheet: "out" channel counts (O column, and L column
# for post-processing rows) changed from 30 to 21, and the derived
# time[us] values (S column) were recalculated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 88 (output1 Conv)
$ws.Range("O88").Value = 21
$ws.Range("S88").Value = 125

# Row 97 (output2 Conv)
$ws.Range("O97").Value = 21
$ws.Range("S97").Value = 226

# Row 106 (output3 Conv)
$ws.Range("O106").Value = 21
$ws.Range("S106").Value = 469

# Row 107 (post_0_0_transpose)
$ws.Range("L107").Value = 21
$ws.Range("O107").Value = 21
$ws.Range("S107").Value = 4

# Row 108 (post_0_1_cast_fp16_fp32)
$ws.Range("L108").Value = 21
$ws.Range("O108").Value = 21
$ws.Range("S108").Value = 3

# Row 109 (post_1_0_transpose)
$ws.Range("L109").Value = 21
$ws.Range("O109").Value = 21
$ws.Range("S109").Value = 14

# Row 110 (post_1_1_cast_fp16_fp32)
$ws.Range("L110").Value = 21
$ws.Range("O110").Value = 21
$ws.Range("S110").Value = 12

# Row 111 (post_2_0_transpose)
$ws.Range("L111").Value = 21
$ws.Range("O111").Value = 21
$ws.Range("S111").Value = 54

# Row 112 (post_2_1_cast_fp16_fp32)
$ws.Range("L112").Value = 21
$ws.Range("O112").Value = 21
$ws.Range("S112").Value = 45

# Row 113 (Total time)
$ws.Range("S113").Value = 188888
